$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate the Id (A), Aktivitet (M), Ost (Q) and Nord (R) values among rows 6, 7 and 8.
# Row 6 <- old Row 8 data, Row 7 <- old Row 6 data, Row 8 <- old Row 7 data.

$ws.Range("A6").Value = 107258607
$ws.Range("M6").Value = ""
$ws.Range("Q6").Value = 404755.5111078721
$ws.Range("R6").Value = 7063764.822795196

$ws.Range("A7").Value = 107258628
$ws.Range("M7").Value = ""
$ws.Range("Q7").Value = 404588.0690095468
$ws.Range("R7").Value = 7064520.029476826

$ws.Range("A8").Value = 107258608
$ws.Range("M8").Value = "färska spår"
$ws.Range("Q8").Value = 404465.3900776547
$ws.Range("R8").Value = 7064504.653031247
